# Generate Report for Handback
#
# This localization-status report is updated to reflect that the handback
# (translated files coming back in sync with en-US) has happened:
#   - "Status" goes from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the Overview sheet and on each language sheet (zh-cn, de-de).
#   - Each language sheet grows two new columns worth of data for row 2/3:
#     F ("Latest Target File") and G ("Latest Handback File"), each populated
#     with a hyperlinked file name (styled like the existing hyperlink cells).
#   - The "Latest Handback DateTime" (column H) is stamped with a real
#     timestamp instead of the zero-date placeholder - a different timestamp
#     per language, since zh-cn and de-de completed their handback at
#     different times.

$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Status column updates (shared text across Overview + each language sheet)
# ---------------------------------------------------------------------------
$ov.Range("B2:C3").Value = $statusText
$zh.Range("C2:C3").Value = $statusText
$de.Range("C2:C3").Value = $statusText

# ---------------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) - real timestamps now, per language
# ---------------------------------------------------------------------------
$zh.Range("H2:H3").Value = "2016-03-24 00:34:20"
$de.Range("H2:H3").Value = "2016-03-24 00:34:29"

# ---------------------------------------------------------------------------
# 3. New columns F (Latest Target File) / G (Latest Handback File) for rows 2 & 3
#    on each language sheet - hyperlinked, styled like the existing link cells.
# ---------------------------------------------------------------------------
function Add-LinkCell($sheet, $cellRef, $text, $url) {
    $range = $sheet.Range($cellRef)
    $range.Value = $text
    $sheet.Hyperlinks.Add($range, $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $text) | Out-Null
    # Match the look of the workbook's other hyperlink cells (underlined,
    # Excel's standard hyperlink blue).
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = $true
    $range.Font.Color = 15570276
}

# zh-cn: reuse the same targets as the existing "a.md" / handoff-xlf links
$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4f290c9cf38377daca6b54f02a1894c25cf4fa7b/e2e/a.md"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/53ed57e98cac4f9a90be735eca34ca616f31520d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

Add-LinkCell $zh "F2" "a.md" $zhMdUrl
Add-LinkCell $zh "G2" $zhXlfName $zhXlfUrl
Add-LinkCell $zh "F3" "a.md" $zhMdUrl
Add-LinkCell $zh "G3" $zhXlfName $zhXlfUrl

# de-de: reuse the same targets as the existing "a.md" / handoff-xlf links
$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4f290c9cf38377daca6b54f02a1894c25cf4fa7b/e2e/a.md"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8166dee6ad3f9907bcd72c65ff6a280c26393b89/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

Add-LinkCell $de "F2" "a.md" $deMdUrl
Add-LinkCell $de "G2" $deXlfName $deXlfUrl
Add-LinkCell $de "F3" "a.md" $deMdUrl
Add-LinkCell $de "G3" $deXlfName $deXlfUrl

Write-Output "Handback report generated."
